$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "Sheet1" summary sheet FIRST (as the leftmost tab).
#    Do this before grabbing references to the other sheets, since
#    inserting a sheet shifts sheet positions/identities.
# ------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Add()
$wsSummary.Name = "Sheet1"

# Now (re-)fetch the existing vendor sheets by name.
$wsAmazon   = $wb.Worksheets.Item("Amazon")
$wsWalmart  = $wb.Worksheets.Item("Walmart")
$wsGamestop = $wb.Worksheets.Item("Gamestop")

# ------------------------------------------------------------------
# 2. Rename a few products (text edits on existing vendor sheets).
# ------------------------------------------------------------------
$wsWalmart.Range("A2").Value = "Hoodie"          # was "Champion Hoodie"
$wsGamestop.Range("A3").Value = "Videogame"      # was "Far Cry 6"
$wsGamestop.Range("A4").Value = "Controller"     # was "Oculus "

# ------------------------------------------------------------------
# 3. Update Gamestop's price/stock numbers for the renamed items.
# ------------------------------------------------------------------
$wsGamestop.Range("B4").Value = 59.99
$wsGamestop.Range("C3").Value = 33
$wsGamestop.Range("C4").Value = 29

# ------------------------------------------------------------------
# 4. Populate the new "Sheet1" with a combined summary table.
# ------------------------------------------------------------------
$wsSummary.Range("A1").Value = "Products:"
$wsSummary.Range("B1").Value = "Prices:"
$wsSummary.Range("C1").Value = "Stock:"

$data = @(
    @("Sunglasses",       15.99,  57),
    @("IPad",             329.99, 25),
    @("Leaf Blower",      169.99, 11),
    @("Hoodie",            21.99, 16),
    @("Basketball Hoop",  169.99,  3),
    @("Blender",           74.99,  0),
    @("PS5",              499.99,  0),
    @("Videogame",         59.99, 33),
    @("Controller",        59.99, 29)
)

$r = 2
foreach ($row in $data) {
    $wsSummary.Cells.Item($r, 1).Value = $row[0]
    $wsSummary.Cells.Item($r, 2).Value = $row[1]
    $wsSummary.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Header row A1:I1 gets a left-aligned style (matches added cellXf).
$wsSummary.Range("A1:I1").HorizontalAlignment = -4131

# ------------------------------------------------------------------
# 5. Restore per-sheet selections / active sheet.
# ------------------------------------------------------------------
$wsAmazon.Activate()
$wsAmazon.Range("B19").Select()

$wsWalmart.Activate()
$wsWalmart.Range("B2").Select()

$wsGamestop.Activate()
$wsGamestop.Range("A5").Select()

$wsSummary.Activate()
$wsSummary.Range("D5").Select()
